# Add new poll and betting data.
# A new "Latest Morgan" poll is added at row 6, pushing the previous
# "Latest Morgan" -> "Second Morgan" (row 7) and the previous
# "Second Morgan" -> "Third Morgan" (row 8). The old "Third Morgan" data
# is discarded. Rows 12/13 (Essential 3-poll averages) are also refreshed
# with updated raw figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# Row 8 <- old Row 7 data (Third Morgan, shifted down)
$ws.Range("B8").Value = 56.5
$ws.Range("C8").Value = 56
$ws.Range("D8").Value = 63.5
$ws.Range("E8").Value = 48.5
$ws.Range("F8").Value = 52
$ws.Range("G8").Value = 52.5

# Row 7 <- old Row 6 data (Second Morgan, shifted down)
$ws.Range("B7").Value = 56
$ws.Range("C7").Value = 56.5
$ws.Range("D7").Value = 60
$ws.Range("E7").Value = 48
$ws.Range("F7").Value = 53
$ws.Range("G7").Value = 54.5

# Row 6 <- brand new poll data (Latest Morgan)
$ws.Range("B6").Value = 58
$ws.Range("C6").Value = 57.5
$ws.Range("D6").Value = 64
$ws.Range("E6").Value = 45.5
$ws.Range("F6").Value = 59
$ws.Range("G6").Value = 60.5

# Row 12 (Latest Essential (3 avg.) ->) updated figures
$ws.Range("B12").Value = 52.329749103942646
$ws.Range("C12").Value = 49.820788530465954
$ws.Range("D12").Value = 51.798561151079141
$ws.Range("F12").Value = 54.838709677419352
$ws.Range("G12").Value = 54.255319148936174

# Row 13 (Second Essential (3 avg.) ->) updated figures
$ws.Range("B13").Value = 52.158273381294968
$ws.Range("C13").Value = 49.81818181818182
$ws.Range("D13").Value = 53.763440860215056
$ws.Range("F13").Value = 57.446808510638299

# Update the active selection to match the author's final cursor position.
$ws.Range("I16").Select() | Out-Null
